$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2 holds rich (mixed-formatting) text. Rewrite its value and
# reapply per-run formatting: "Bold" stays bold, "Italic" stays italic,
# the separators become plain/normal text.
$a2 = $ws.Range("A2")
$a2.Value = "Bold Italic and plain"

$a2.Characters(1, 4).Font.Bold = $true
$a2.Characters(1, 4).Font.Italic = $false

$a2.Characters(5, 1).Font.Bold = $false
$a2.Characters(5, 1).Font.Italic = $false

$a2.Characters(6, 6).Font.Bold = $false
$a2.Characters(6, 6).Font.Italic = $true

$a2.Characters(12, 10).Font.Bold = $false
$a2.Characters(12, 10).Font.Italic = $false

# B2 / C2 keep their original displayed values ("Normal" / "Bold Only");
# rewrite them too so the shared-strings table is rebuilt the same way
# inline strings would be handled.
$ws.Range("B2").Value = "Normal"
$ws.Range("C2").Value = "Bold Only"

# Move the active selection to B2.
$ws.Range("B2").Select() | Out-Null
